$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous upload's manual formatting (direct font/style on every cell)
# is being cleared back to the workbook default style before the fresh
# Jobsheet 10 data is written in. D2 is left alone here so it keeps its
# existing date-number-format style (only its value changes below).
$ws.Range("A1:E1").Style = "Normal"
$ws.Range("A2:C2").Style = "Normal"
$ws.Range("E2").Style = "Normal"

# Row 1 (headers) keep their text - only the styling above changed.

# Row 2: new stock entry for "Jobsheet 10_Praktikum 1"
$ws.Range("A2").Value = 2
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 45588
$ws.Range("E2").Value = 150

# Cursor was left on S1 when the sheet was saved
$ws.Range("S1").Select() | Out-Null
